$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 69
$ws.Cells.Item($newRow, 1).Value = "'2021/08/18"
$ws.Cells.Item($newRow, 2).Value = 314.6
$ws.Cells.Item($newRow, 3).Value = 319
$ws.Cells.Item($newRow, 4).Value = 0.98
$ws.Cells.Item($newRow, 5).Value = 0.98
